# Updates cryptos list prices/volume percentages (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column ("Price") values are written as literal text: many look like plain
# numbers (or use "." as a thousands separator, e.g. "63.785.68"), and Excel
# would otherwise silently reinterpret/reformat them as numeric values and lose
# the original text (e.g. "8.10" -> 8.1, "0.999" -> 0.999 as a float with no
# trailing digits preserved textually). Forcing the cell to Text format first,
# then clearing the format back (so no stray style sticks to the cell), keeps
# the exact source string while leaving cell formatting untouched.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '63.785.68'
$ws.Range('E2').Value = '  -1.72%  '
Set-TextValue 'D3' '3.048.59'
$ws.Range('E3').Value = '  -1.92%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue 'D5' '555.91'
$ws.Range('E5').Value = '  -0.68%  '
Set-TextValue 'D6' '141.75'
$ws.Range('E6').Value = '  -2.37%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue 'D8' '3.048.95'
$ws.Range('E8').Value = '  -1.75%  '
Set-TextValue 'D9' '0.522'
$ws.Range('E9').Value = '  +4.14%  '
Set-TextValue 'D10' '0.153'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  -11.80%  '
Set-TextValue 'D12' '0.485'
$ws.Range('E12').Value = '  +3.69%  '
Set-TextValue 'D13' '0.0000229'
$ws.Range('E13').Value = '  -1.68%  '
Set-TextValue 'D14' '35.24'
$ws.Range('E14').Value = '  -0.90%  '
Set-TextValue 'D15' '3.542.89'
$ws.Range('E15').Value = '  -1.67%  '
Set-TextValue 'D16' '63.736.18'
$ws.Range('E16').Value = '  -1.85%  '
Set-TextValue 'D17' '3.043.91'
$ws.Range('E17').Value = '  -1.99%  '
Set-TextValue 'D19' '6.76'
$ws.Range('E19').Value = '  -1.35%  '
Set-TextValue 'D20' '483.92'
$ws.Range('E20').Value = '  -0.08%  '
Set-TextValue 'D21' '14.11'
$ws.Range('E21').Value = '  +1.82%  '
Set-TextValue 'D22' '0.682'
$ws.Range('E22').Value = '  +0.03%  '
Set-TextValue 'D23' '14.51'
$ws.Range('E23').Value = '  +7.72%  '
Set-TextValue 'D24' '7.52'
$ws.Range('E24').Value = '  -0.10%  '
Set-TextValue 'D25' '82.51'
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.21%  '
Set-TextValue 'D28' '8.10'
$ws.Range('E28').Value = '  -1.44%  '
Set-TextValue 'D29' '2.03'
$ws.Range('E29').Value = '  -2.04%  '
Set-TextValue 'D30' '0.998'
$ws.Range('E30').Value = '  -0.16%  '
Set-TextValue 'D31' '26.15'
Set-TextValue 'D32' '1.15'
$ws.Range('E32').Value = '  -0.86%  '
Set-TextValue 'D33' '2.44'
$ws.Range('E33').Value = '  -1.43%  '
Set-TextValue 'D34' '5.65'
$ws.Range('E34').Value = '  -1.11%  '
Set-TextValue 'D35' '6.20'
$ws.Range('E35').Value = '  -0.82%  '
Set-TextValue 'D36' '55.10'
$ws.Range('E36').Value = '  +0.12%  '
Set-TextValue 'D37' '0.0407'
$ws.Range('E37').Value = '  -0.92%  '
Set-TextValue 'D38' '441.06'
$ws.Range('E38').Value = '  -6.14%  '
Set-TextValue 'D39' '0.0815'
$ws.Range('E39').Value = '  -2.02%  '
Set-TextValue 'D40' '3.003.29'
$ws.Range('E40').Value = '  -0.80%  '
Set-TextValue 'D41' '2.76'
$ws.Range('E41').Value = '  -5.02%  '
Set-TextValue 'D42' '8.30'
$ws.Range('E42').Value = '  -0.10%  '
Set-TextValue 'D43' '0.115'
$ws.Range('E43').Value = '  -1.20%  '
Set-TextValue 'D44' '0.269'
$ws.Range('E44').Value = '  +3.66%  '
Set-TextValue 'D45' '27.65'
$ws.Range('E45').Value = '  -3.07%  '
Set-TextValue 'D46' '2.23'
$ws.Range('E46').Value = '  +4.69%  '
$ws.Range('E48').Value = '  +0.06%  '
Set-TextValue 'D49' '118.03'
$ws.Range('E49').Value = '  -0.08%  '
Set-TextValue 'D50' '0.0₃0511'
$ws.Range('E50').Value = '  -1.43%  '
Set-TextValue 'D51' '2.09'
$ws.Range('E51').Value = '  +0.06%  '
